# Fill in the "Number of orders_sf" (column C) and "differ_Number of orders"
# (column D) values for the Sales Organization rows that were previously
# blank/zero because the SF result file had not yet been downloaded.
# For each of these rows B (Number of orders_sap) is 0, so D = |B - C| = C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    27  = 238
    29  = 1034
    32  = 186
    35  = 4417
    39  = 32
    40  = 61
    42  = 8233
    44  = 1806
    46  = 2812
    51  = 1956
    54  = 133
    64  = 111
    71  = 946
    72  = 7764
    74  = 61
    79  = 138
    84  = 3143
    86  = 1285
    88  = 32
    96  = 87
    99  = 1085
    101 = 36
    104 = 12
    107 = 4888
    111 = 248415
    112 = 65
    118 = 106
}

foreach ($row in $updates.Keys) {
    $value = $updates[$row]
    $ws.Cells.Item($row, 3).Value = $value
    $ws.Cells.Item($row, 4).Value = $value
}
